$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---
# Column B (Week_Start_Date) is stored as text, not a true date, so force
# text storage (otherwise Excel auto-converts "yyyy-mm-dd" looking values
# into a date serial number).
$newWeekStart = @{
  2  = "2024-12-29"
  3  = "2025-01-05"
  4  = "2025-01-12"
  5  = "2025-01-19"
  6  = "2025-01-26"
  7  = "2025-02-02"
  8  = "2025-02-09"
  9  = "2025-02-16"
  10 = "2025-02-23"
  11 = "2025-03-02"
  12 = "2025-03-09"
  13 = "2025-03-16"
  14 = "2025-03-23"
  15 = "2025-03-30"
  16 = "2025-04-06"
  17 = "2025-04-13"
}

# Column D (MyForecast) updated values (only rows that actually change).
$newForecast = @{
  2  = 5
  4  = 2
  5  = 3
  6  = 5
  7  = 3
  8  = 6
  10 = 3
  12 = 6
  13 = 8
  14 = 6
  16 = 6
}

for ($r = 2; $r -le 17; $r++) {
  $cellB = $ws.Cells.Item($r, 2)
  $cellB.NumberFormat = "@"
  $cellB.Value = $newWeekStart[$r]

  if ($newForecast.ContainsKey($r)) {
    $ws.Cells.Item($r, 4).Value = $newForecast[$r]
  }

  # Column J (is_holiday_week) is cleared out (becomes an empty/blank cell).
  $ws.Cells.Item($r, 10).ClearContents()
}

# --- Summary sheet ---
# All "Value" entries on this sheet are stored as plain text, so force text
# storage before assigning (otherwise numeric-looking / date-looking strings
# get auto-converted by Excel).
$summaryUpdates = @{
  9  = "75"
  10 = "34"
  12 = "8"
  13 = "2025-03-16"
  14 = "2"
  15 = "2025-01-12"
}

foreach ($r in $summaryUpdates.Keys) {
  $cell = $summary.Cells.Item($r, 2)
  $cell.NumberFormat = "@"
  $cell.Value = $summaryUpdates[$r]
}
